$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe in the literal forces Excel to store numeric-looking
# price strings (e.g. "37.252.56", "247.50", "0.628") as literal text instead
# of auto-converting/rounding them as numbers.

$ws.Range("D2").Value = '''37.252.56'
$ws.Range("E2").Value = '  +1.54%  '

$ws.Range("D3").Value = '''2.017.99'

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''247.50'
$ws.Range("E5").Value = '  +1.02%  '

$ws.Range("D6").Value = '''0.628'
$ws.Range("E6").Value = '  +1.45%  '

$ws.Range("D7").Value = '''60.11'
$ws.Range("E7").Value = '  -2.24%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '''0.389'
$ws.Range("E9").Value = '  +2.90%  '

$ws.Range("E10").Value = '  +1.17%  '

$ws.Range("E11").Value = '  +1.12%  '

$ws.Range("D12").Value = '''15.15'
$ws.Range("E12").Value = '  +5.73%  '

$ws.Range("D13").Value = '''2.312.59'
$ws.Range("E13").Value = '  +3.44%  '

$ws.Range("D14").Value = '''0.852'
$ws.Range("E14").Value = '  +1.30%  '

$ws.Range("E15").Value = '  +1.07%  '

$ws.Range("D16").Value = '''5.56'
$ws.Range("E16").Value = '  +4.52%  '

$ws.Range("D17").Value = '''2.019.22'
$ws.Range("E17").Value = '  +3.04%  '

$ws.Range("D18").Value = '''37.153.98'
$ws.Range("E18").Value = '  +1.25%  '

$ws.Range("D19").Value = '''70.58'
$ws.Range("E19").Value = '  +0.75%  '

$ws.Range("D20").Value = '''0.0₃0870'
$ws.Range("E20").Value = '  +1.51%  '

$ws.Range("D21").Value = '''5.24'
$ws.Range("E21").Value = '  +2.56%  '

$ws.Range("D22").Value = '''231.18'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").Value = '''2.51'
$ws.Range("E24").Value = '  +1.32%  '

$ws.Range("E25").Value = '  +0.43%  '

$ws.Range("D26").Value = '''9.47'
$ws.Range("E26").Value = '  +2.37%  '

$ws.Range("E27").Value = '  +2.14%  '

$ws.Range("D28").Value = '''0.139'
$ws.Range("E28").Value = '  -3.38%  '

$ws.Range("D29").Value = '''19.78'
$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("E30").Value = '  +11.84%  '

$ws.Range("E31").Value = '  +1.52%  '

$ws.Range("D32").Value = '''4.83'
$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("D33").Value = '''0.0664'
$ws.Range("E33").Value = '  +6.82%  '

$ws.Range("D34").Value = '''4.50'
$ws.Range("E34").Value = '  +0.31%  '

$ws.Range("E35").Value = '  +8.14%  '

$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").Value = '''3.50'
$ws.Range("E36").Value = '  -0.96%  '

$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.09%  '

$ws.Range("E38").Value = '  +2.10%  '

$ws.Range("D39").Value = '''5.42'
$ws.Range("E39").Value = '  -3.00%  '

$ws.Range("D40").Value = '''0.0987'
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("E41").Value = '  +0.68%  '

$ws.Range("E42").Value = '  +1.03%  '

$ws.Range("E43").Value = '  +1.30%  '

$ws.Range("D44").Value = '''16.74'
$ws.Range("E44").Value = '  +3.13%  '

$ws.Range("D45").Value = '''92.40'
$ws.Range("E45").Value = '  +3.70%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''1.07'
$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '''1.386.59'
$ws.Range("E47").Value = '  +1.22%  '

$ws.Range("D48").Value = '''7.55'
$ws.Range("E48").Value = '  +5.02%  '

$ws.Range("D49").Value = '''2.12'
$ws.Range("E49").Value = '  +13.62%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '''47.04'
$ws.Range("E50").Value = '  +5.37%  '

$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").Value = '''2.85'
$ws.Range("E51").Value = '  +0.20%  '
